# Insert a new price record as row 116 on the active sheet, pushing the
# existing rows 116-121 down to 117-122 (dimension grows from A1:T121 to
# A1:T122).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 116 (shifts rows 116:121 -> 117:122).
$ws.Rows.Item(116).Insert()

# Populate the new row 116 with the new weekly record.
$row = 116
$ws.Cells.Item($row, 1).Value2 = 10
$ws.Cells.Item($row, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value2 = "La Araucanía"
$ws.Cells.Item($row, 4).Value2 = 44615
$ws.Cells.Item($row, 5).Value2 = 9
$ws.Cells.Item($row, 6).Value2 = "Fruta"
$ws.Cells.Item($row, 7).Value2 = 100104
$ws.Cells.Item($row, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item($row, 9).Value2 = 100104003
$ws.Cells.Item($row, 10).Value2 = "Membrillo"
$ws.Cells.Item($row, 11).Value2 = "Champion"
$ws.Cells.Item($row, 12).Value2 = "Primera"
$ws.Cells.Item($row, 13).Value2 = 50
$ws.Cells.Item($row, 14).Value2 = 16000
$ws.Cells.Item($row, 15).Value2 = 16000
$ws.Cells.Item($row, 16).Value2 = 16000
$ws.Cells.Item($row, 17).Value2 = "`$/bandeja 18 kilos granel"
$ws.Cells.Item($row, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value2 = 889
$ws.Cells.Item($row, 20).Value2 = 18
